# Apply the "DOM Changes in PROD" edit:
# Testdata sheet (B11, row "ContactUSHeader") text is corrected from
# "Contact Kaman Industrial Technologies" to "Contact Us".

$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("TC03_Verify_BLP_Sol_Contact_etc")
$wsData = $wb.Worksheets.Item("Testdata")

$wsData.Range("B11").Value = "Contact Us"

# Match the resulting view/selection state recorded in the workbook:
# the main sheet keeps a plain selection on C13 (no longer the active tab),
# and the Testdata sheet becomes the active tab with B11 selected.
$wsMain.Activate()
$wsMain.Range("C13").Select()

$wsData.Activate()
$wsData.Range("B11").Select()
